$wb = $excel.ActiveWorkbook

# Sheet 1: "Metadata"
$metadata = $wb.Worksheets.Item("Metadata")
$metadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/ada-tooth-quadrant"
$metadata.Range("B3").Value = "8.0.0"
$metadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$metadata.Range("B9").Value = "LinuxForHealth Team"

# Sheet 2: "Include from American Dental "
$codesystem = $wb.Worksheets.Item("Include from American Dental ")
$codesystem.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/ada-tooth-quadrant"
